{"js": "// Office.js (Word JavaScript API) script.\n// Inserts a blank line, a horizontal-rule-style paragraph (bottom border),\n// another blank line, and a full new \"FOCACCIA\" recipe section (heading,\n// ingredient list, and step-by-step instructions) right after the\n// paragraph that ends the FLAN CASERO recipe (\"... As\u00ed de f\u00e1cil.\"),\n// leaving the document's existing trailing empty paragraph untouched.\n\n// Flat-OPC wrapped OOXML fragment for the new paragraphs. Word.Range.insertOoxml\n// requires a full <pkg:package> wrapper (not a bare fragment).\nconst NEW_PARAGRAPHS_OOXML = \"<w:p><w:pPr><w:rPr><w:lang w:val=\\\"es-ES\\\"/></w:rPr></w:pPr></w:p><w:p><w:pPr><w:pBdr><w:bottom w:val=\\\"single\\\" w:sz=\\\"6\\\" w:space=\\\"1\\\" w:color=\\\"auto\\\"/></w:pBdr><w:rPr><w:lang w:val=\\\"es-ES\\\"/></w:rPr></w:pPr></w:p><w:p><w:pPr><w:rPr><w:lang w:val=\\\"es-ES\\\"/></w:rPr></w:pPr></w:p><w:p><w:pPr><w:rPr><w:b/><w:bCs/></w:rPr></w:pPr><w:r><w:rPr><w:b/><w:bCs/></w:rPr><w:t>FOCACCIA</w:t></w:r></w:p><w:p/><w:p><w:r><w:t>Ingredientes</w:t></w:r></w:p><w:p><w:r><w:t>400 g de harina 0000</w:t></w:r></w:p><w:p><w:r><w:t>15 g de levadura de cerveza fresca</w:t></w:r></w:p><w:p><w:r><w:t>280 ml de agua tibia</w:t></w:r></w:p><w:p><w:r><w:t>1\\u20442 vaso de agua para el final</w:t></w:r></w:p><w:p><w:r><w:t>3 cebollas cortadas en juliana</w:t></w:r></w:p><w:p><w:r><w:t>1 diente de ajo grande</w:t></w:r></w:p><w:p><w:r><w:t>100 g de aceitunas (verdes o negras) descarozadas y aplastadas</w:t></w:r></w:p><w:p><w:r><w:t>4 pocillos de aceite de oliva</w:t></w:r></w:p><w:p><w:r><w:t>hojas de romero fresco</w:t></w:r></w:p><w:p><w:r><w:t>sal fina y sal gruesa</w:t></w:r></w:p><w:p/><w:p><w:r><w:t>Preparaci\\u00f3n</w:t></w:r></w:p><w:p><w:r><w:t xml:space=\\\"preserve\\\">\\u00b7 </w:t></w:r><w:r><w:t>Cocinar la cebolla junto con un pocillo de aceite de oliva a fuego lento en una olla tapada. Salpimentar.</w:t></w:r></w:p><w:p><w:r><w:t xml:space=\\\"preserve\\\">\\u00b7 </w:t></w:r><w:r><w:t>Cuando est\\u00e9 cocida, retirar del calor y reservar.</w:t></w:r></w:p><w:p><w:r><w:t xml:space=\\\"preserve\\\">\\u00b7 </w:t></w:r><w:r><w:t>Licuar el ajo con el agua tibia.</w:t></w:r></w:p><w:p><w:r><w:t xml:space=\\\"preserve\\\">\\u00b7 </w:t></w:r><w:r><w:t>Deshacer la levadura e incorporar al ajo junto con un pocillo de aceite. Reservar.</w:t></w:r></w:p><w:p><w:r><w:t xml:space=\\\"preserve\\\">\\u00b7 </w:t></w:r><w:r><w:t>En otro recipiente, mezclar la harina con una cucharada de sal y 2 pocillos de aceite.</w:t></w:r></w:p><w:p><w:r><w:t xml:space=\\\"preserve\\\">\\u00b7 </w:t></w:r><w:r><w:t>Incorporar las dos preparaciones anteriores.</w:t></w:r></w:p><w:p><w:r><w:t xml:space=\\\"preserve\\\">\\u00b7 </w:t></w:r><w:r><w:t>Mezclar hasta que se forme una masa homog\\u00e9nea.</w:t></w:r></w:p><w:p><w:r><w:t xml:space=\\\"preserve\\\">\\u00b7 </w:t></w:r><w:r><w:t xml:space=\\\"preserve\\\">Aceitar un molde rectangular y colocar la masa que debe estar algo </w:t></w:r><w:proofErr w:type=\\\"spellStart\\\"/><w:r><w:t>h\\u00fameda</w:t></w:r><w:proofErr w:type=\\\"spellEnd\\\"/><w:r><w:t>.</w:t></w:r></w:p><w:p><w:r><w:t xml:space=\\\"preserve\\\">\\u00b7 </w:t></w:r><w:r><w:t xml:space=\\\"preserve\\\">Dejar levar tapada hasta que crezca un par de </w:t></w:r><w:proofErr w:type=\\\"spellStart\\\"/><w:r><w:t>cent\\u00edmetros</w:t></w:r><w:proofErr w:type=\\\"spellEnd\\\"/><w:r><w:t>.</w:t></w:r></w:p><w:p><w:r><w:t xml:space=\\\"preserve\\\">\\u00b7 </w:t></w:r><w:r><w:t>Desparramar sobre la masa las aceitunas y las hojas de romero, hundi\\u00e9ndolas con las yemas de los dedos.</w:t></w:r></w:p><w:p><w:r><w:t xml:space=\\\"preserve\\\">\\u00b7 </w:t></w:r><w:r><w:t xml:space=\\\"preserve\\\">Dejar leudar nuevamente no </w:t></w:r><w:proofErr w:type=\\\"spellStart\\\"/><w:r><w:t>m\\u00e1s</w:t></w:r><w:proofErr w:type=\\\"spellEnd\\\"/><w:r><w:t xml:space=\\\"preserve\\\"> de 15 minutos.</w:t></w:r></w:p><w:p><w:r><w:t xml:space=\\\"preserve\\\">\\u00b7 </w:t></w:r><w:r><w:t>Espolvorear ligeramente con la sal gruesa.</w:t></w:r></w:p><w:p><w:r><w:t xml:space=\\\"preserve\\\">\\u00b7 </w:t></w:r><w:r><w:t xml:space=\\\"preserve\\\">Rociar con agua para que quede </w:t></w:r><w:proofErr w:type=\\\"spellStart\\\"/><w:r><w:t>h\\u00famedo</w:t></w:r><w:proofErr w:type=\\\"spellEnd\\\"/><w:r><w:t>.</w:t></w:r></w:p><w:p><w:r><w:t xml:space=\\\"preserve\\\">\\u00b7 </w:t></w:r><w:r><w:t>Cocinar en horno bien caliente entre 20 y 25 minutos.</w:t></w:r></w:p><w:p><w:r><w:t xml:space=\\\"preserve\\\">\\u00b7 </w:t></w:r><w:r><w:t>Retirar la </w:t></w:r><w:proofErr w:type=\\\"spellStart\\\"/><w:r><w:rPr><w:b/><w:bCs/></w:rPr><w:t>focaccia</w:t></w:r><w:proofErr w:type=\\\"spellEnd\\\"/><w:r><w:rPr><w:b/><w:bCs/></w:rPr><w:t xml:space=\\\"preserve\\\"> con oliva,</w:t></w:r><w:r><w:t> dejar enfriar y cortar.</w:t></w:r></w:p><w:p/>\";\n\nconst FLAT_OPC_TEMPLATE = ooxmlBody => `<?xml version=\"1.0\" encoding=\"UTF-8\" standalone=\"yes\"?>\n<pkg:package xmlns:pkg=\"http://schemas.microsoft.com/office/2006/xmlPackage\">\n  <pkg:part pkg:name=\"/_rels/.rels\" pkg:contentType=\"application/vnd.openxmlformats-package.relationships+xml\">\n    <pkg:xmlData>\n      <Relationships xmlns=\"http://schemas.openxmlformats.org/package/2006/relationships\">\n        <Relationship Id=\"rId1\" Type=\"http://schemas.openxmlformats.org/officeDocument/2006/relationships/officeDocument\" Target=\"word/document.xml\"/>\n      </Relationships>\n    </pkg:xmlData>\n  </pkg:part>\n  <pkg:part pkg:name=\"/word/document.xml\" pkg:contentType=\"application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml\">\n    <pkg:xmlData>\n      <w:document xmlns:w=\"http://schemas.openxmlformats.org/wordprocessingml/2006/main\">\n        <w:body>${ooxmlBody}</w:body>\n      </w:document>\n    </pkg:xmlData>\n  </pkg:part>\n</pkg:package>`;\n\nconst body = context.document.body;\nconst paragraphs = body.paragraphs;\nparagraphs.load(\"items/text\");\nawait context.sync();\n\n// Locate the paragraph that ends with \"As\u00ed de f\u00e1cil.\" (the last line of the\n// FLAN CASERO recipe) -- the new content is inserted right after it, before\n// the document's final (trailing) paragraph.\nlet anchor = null;\nfor (let i = 0; i < paragraphs.items.length; i++) {\n  if (paragraphs.items[i].text.indexOf(\"As\u00ed de f\u00e1cil.\") !== -1) {\n    anchor = paragraphs.items[i];\n    break;\n  }\n}\nif (!anchor) {\n  throw new Error('Could not find the anchor paragraph (\"As\u00ed de f\u00e1cil.\")');\n}\n\nconst insertionPoint = anchor.getRange(\"End\");\ninsertionPoint.insertOoxml(FLAT_OPC_TEMPLATE(NEW_PARAGRAPHS_OOXML), \"After\");\nawait context.sync();\n", "ps1": "# Word COM interop (PowerShell-style) script.\n# Inserts a blank line, a horizontal-rule-style paragraph (bottom border),\n# another blank line, and a full new \"FOCACCIA\" recipe section (heading,\n# ingredient list, and step-by-step instructions) right after the\n# paragraph that ends the FLAN CASERO recipe (\"... Asi de facil.\"),\n# leaving the document's existing trailing empty paragraph untouched.\n\n$d = $word.ActiveDocument\n\n# Flat-OPC wrapped OOXML fragment for the new paragraphs. Range.InsertXML\n# requires a full <pkg:package> wrapper (not a bare fragment).\n$newParagraphsOoxml = '<w:p><w:pPr><w:rPr><w:lang w:val=\"es-ES\"/></w:rPr></w:pPr></w:p><w:p><w:pPr><w:pBdr><w:bottom w:val=\"single\" w:sz=\"6\" w:space=\"1\" w:color=\"auto\"/></w:pBdr><w:rPr><w:lang w:val=\"es-ES\"/></w:rPr></w:pPr></w:p><w:p><w:pPr><w:rPr><w:lang w:val=\"es-ES\"/></w:rPr></w:pPr></w:p><w:p><w:pPr><w:rPr><w:b/><w:bCs/></w:rPr></w:pPr><w:r><w:rPr><w:b/><w:bCs/></w:rPr><w:t>FOCACCIA</w:t></w:r></w:p><w:p/><w:p><w:r><w:t>Ingredientes</w:t></w:r></w:p><w:p><w:r><w:t>400 g de harina 0000</w:t></w:r></w:p><w:p><w:r><w:t>15 g de levadura de cerveza fresca</w:t></w:r></w:p><w:p><w:r><w:t>280 ml de agua tibia</w:t></w:r></w:p><w:p><w:r><w:t>1\u20442 vaso de agua para el final</w:t></w:r></w:p><w:p><w:r><w:t>3 cebollas cortadas en juliana</w:t></w:r></w:p><w:p><w:r><w:t>1 diente de ajo grande</w:t></w:r></w:p><w:p><w:r><w:t>100 g de aceitunas (verdes o negras) descarozadas y aplastadas</w:t></w:r></w:p><w:p><w:r><w:t>4 pocillos de aceite de oliva</w:t></w:r></w:p><w:p><w:r><w:t>hojas de romero fresco</w:t></w:r></w:p><w:p><w:r><w:t>sal fina y sal gruesa</w:t></w:r></w:p><w:p/><w:p><w:r><w:t>Preparaci\u00f3n</w:t></w:r></w:p><w:p><w:r><w:t xml:space=\"preserve\">\u00b7 </w:t></w:r><w:r><w:t>Cocinar la cebolla junto con un pocillo de aceite de oliva a fuego lento en una olla tapada. Salpimentar.</w:t></w:r></w:p><w:p><w:r><w:t xml:space=\"preserve\">\u00b7 </w:t></w:r><w:r><w:t>Cuando est\u00e9 cocida, retirar del calor y reservar.</w:t></w:r></w:p><w:p><w:r><w:t xml:space=\"preserve\">\u00b7 </w:t></w:r><w:r><w:t>Licuar el ajo con el agua tibia.</w:t></w:r></w:p><w:p><w:r><w:t xml:space=\"preserve\">\u00b7 </w:t></w:r><w:r><w:t>Deshacer la levadura e incorporar al ajo junto con un pocillo de aceite. Reservar.</w:t></w:r></w:p><w:p><w:r><w:t xml:space=\"preserve\">\u00b7 </w:t></w:r><w:r><w:t>En otro recipiente, mezclar la harina con una cucharada de sal y 2 pocillos de aceite.</w:t></w:r></w:p><w:p><w:r><w:t xml:space=\"preserve\">\u00b7 </w:t></w:r><w:r><w:t>Incorporar las dos preparaciones anteriores.</w:t></w:r></w:p><w:p><w:r><w:t xml:space=\"preserve\">\u00b7 </w:t></w:r><w:r><w:t>Mezclar hasta que se forme una masa homog\u00e9nea.</w:t></w:r></w:p><w:p><w:r><w:t xml:space=\"preserve\">\u00b7 </w:t></w:r><w:r><w:t xml:space=\"preserve\">Aceitar un molde rectangular y colocar la masa que debe estar algo </w:t></w:r><w:proofErr w:type=\"spellStart\"/><w:r><w:t>h\u00fameda</w:t></w:r><w:proofErr w:type=\"spellEnd\"/><w:r><w:t>.</w:t></w:r></w:p><w:p><w:r><w:t xml:space=\"preserve\">\u00b7 </w:t></w:r><w:r><w:t xml:space=\"preserve\">Dejar levar tapada hasta que crezca un par de </w:t></w:r><w:proofErr w:type=\"spellStart\"/><w:r><w:t>cent\u00edmetros</w:t></w:r><w:proofErr w:type=\"spellEnd\"/><w:r><w:t>.</w:t></w:r></w:p><w:p><w:r><w:t xml:space=\"preserve\">\u00b7 </w:t></w:r><w:r><w:t>Desparramar sobre la masa las aceitunas y las hojas de romero, hundi\u00e9ndolas con las yemas de los dedos.</w:t></w:r></w:p><w:p><w:r><w:t xml:space=\"preserve\">\u00b7 </w:t></w:r><w:r><w:t xml:space=\"preserve\">Dejar leudar nuevamente no </w:t></w:r><w:proofErr w:type=\"spellStart\"/><w:r><w:t>m\u00e1s</w:t></w:r><w:proofErr w:type=\"spellEnd\"/><w:r><w:t xml:space=\"preserve\"> de 15 minutos.</w:t></w:r></w:p><w:p><w:r><w:t xml:space=\"preserve\">\u00b7 </w:t></w:r><w:r><w:t>Espolvorear ligeramente con la sal gruesa.</w:t></w:r></w:p><w:p><w:r><w:t xml:space=\"preserve\">\u00b7 </w:t></w:r><w:r><w:t xml:space=\"preserve\">Rociar con agua para que quede </w:t></w:r><w:proofErr w:type=\"spellStart\"/><w:r><w:t>h\u00famedo</w:t></w:r><w:proofErr w:type=\"spellEnd\"/><w:r><w:t>.</w:t></w:r></w:p><w:p><w:r><w:t xml:space=\"preserve\">\u00b7 </w:t></w:r><w:r><w:t>Cocinar en horno bien caliente entre 20 y 25 minutos.</w:t></w:r></w:p><w:p><w:r><w:t xml:space=\"preserve\">\u00b7 </w:t></w:r><w:r><w:t>Retirar la </w:t></w:r><w:proofErr w:type=\"spellStart\"/><w:r><w:rPr><w:b/><w:bCs/></w:rPr><w:t>focaccia</w:t></w:r><w:proofErr w:type=\"spellEnd\"/><w:r><w:rPr><w:b/><w:bCs/></w:rPr><w:t xml:space=\"preserve\"> con oliva,</w:t></w:r><w:r><w:t> dejar enfriar y cortar.</w:t></w:r></w:p><w:p/>'\n\n$flatOpcXml = '<?xml version=\"1.0\" encoding=\"UTF-8\" standalone=\"yes\"?>' +\n  '<pkg:package xmlns:pkg=\"http://schemas.microsoft.com/office/2006/xmlPackage\">' +\n  '<pkg:part pkg:name=\"/_rels/.rels\" pkg:contentType=\"application/vnd.openxmlformats-package.relationships+xml\">' +\n  '<pkg:xmlData>' +\n  '<Relationships xmlns=\"http://schemas.openxmlformats.org/package/2006/relationships\">' +\n  '<Relationship Id=\"rId1\" Type=\"http://schemas.openxmlformats.org/officeDocument/2006/relationships/officeDocument\" Target=\"word/document.xml\"/>' +\n  '</Relationships>' +\n  '</pkg:xmlData>' +\n  '</pkg:part>' +\n  '<pkg:part pkg:name=\"/word/document.xml\" pkg:contentType=\"application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml\">' +\n  '<pkg:xmlData>' +\n  '<w:document xmlns:w=\"http://schemas.openxmlformats.org/wordprocessingml/2006/main\">' +\n  \"<w:body>$newParagraphsOoxml</w:body>\" +\n  '</w:document>' +\n  '</pkg:xmlData>' +\n  '</pkg:part>' +\n  '</pkg:package>'\n\n# Locate the paragraph that ends with \"Asi de facil.\" (the last line of the\n# FLAN CASERO recipe) -- the new content is inserted right after it, before\n# the document's final (trailing) paragraph.\n$finder = $d.Content\n$finder.Find.ClearFormatting()\n$finder.Find.Text = \"As\u00ed de f\u00e1cil.\"\n$found = $finder.Find.Execute()\nif (-not $found) {\n    throw 'Could not find the anchor paragraph (\"As\u00ed de f\u00e1cil.\")'\n}\n$anchorParagraph = $finder.Paragraphs(1)\n\n# Collapse a range to the point right after the anchor paragraph's mark\n# (i.e. right before the document's trailing paragraph) and insert there.\n$insertionPoint = $d.Range($anchorParagraph.Range.End, $anchorParagraph.Range.End)\n[void]$insertionPoint.InsertXML($flatOpcXml)\n"}
